$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.545.10"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "3.159.70"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.93"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.07"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "3.151.05"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.472"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.65"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").Value = "3.683.31"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").Value = "64.520.90"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "3.158.60"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.68"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.72"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.716"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.05"
$ws.Range("E24").Value = "  +5.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.52"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.77"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.46"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("E29").Value = "  +3.31%  "
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("E31").Value = "  -7.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.17"
$ws.Range("E32").Value = "  +3.87%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.69"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.08"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "53.02"
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.05"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "444.31"
$ws.Range("E40").Value = "  -2.81%  "
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.31"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "2.882.75"
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("E46").Value = "  +3.06%  "
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.17"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("E51").Value = "  +2.41%  "
